$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.929.12"
$ws.Range("E2").Value = "  +4.34%  "

$ws.Range("D3").Value = "3.036.02"
$ws.Range("E3").Value = "  +3.93%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "566.35"
$ws.Range("E5").Value = "  +3.43%  "

$ws.Range("D6").Value = "141.60"
$ws.Range("E6").Value = "  +8.87%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "0.523"
$ws.Range("E8").Value = "  +1.93%  "

$ws.Range("D9").Value = "3.024.73"
$ws.Range("E9").Value = "  +3.83%  "

$ws.Range("E10").Value = "  +7.25%  "

$ws.Range("D11").Value = "5.27"
$ws.Range("E11").Value = "  +11.20%  "

$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  +3.97%  "

$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").Value = "  +6.10%  "

$ws.Range("D14").Value = "34.26"
$ws.Range("E14").Value = "  +4.51%  "

$ws.Range("E15").Value = "  +1.80%  "

$ws.Range("D16").Value = "3.516.98"
$ws.Range("E16").Value = "  +3.37%  "

$ws.Range("D17").Value = "7.31"
$ws.Range("E17").Value = "  +6.90%  "

$ws.Range("D18").Value = "3.009.11"
$ws.Range("E18").Value = "  +3.20%  "

$ws.Range("D19").Value = "59.821.60"
$ws.Range("E19").Value = "  +4.11%  "

$ws.Range("D20").Value = "439.73"
$ws.Range("E20").Value = "  +5.77%  "

$ws.Range("D21").Value = "13.73"
$ws.Range("E21").Value = "  +4.85%  "

$ws.Range("D22").Value = "0.731"
$ws.Range("E22").Value = "  +7.22%  "

$ws.Range("D23").Value = "7.19"
$ws.Range("E23").Value = "  +3.82%  "

$ws.Range("D24").Value = "13.37"
$ws.Range("E24").Value = "  +2.64%  "

$ws.Range("D25").Value = "81.08"
$ws.Range("E25").Value = "  +1.93%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").Value = "2.27"
$ws.Range("E27").Value = "  +14.61%  "

$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").Value = "2.57"
$ws.Range("E29").Value = "  +4.45%  "

$ws.Range("D30").Value = "7.93"
$ws.Range("E30").Value = "  +6.71%  "

$ws.Range("D31").Value = "26.15"
$ws.Range("E31").Value = "  +4.01%  "

$ws.Range("D32").Value = "6.30"
$ws.Range("E32").Value = "  +5.90%  "

$ws.Range("E33").Value = "  +6.16%  "

$ws.Range("D34").Value = "0.0₃0800"
$ws.Range("E34").Value = "  +18.41%  "

$ws.Range("D37").Value = "2.14"
$ws.Range("E37").Value = "  +3.80%  "

$ws.Range("D38").Value = "49.19"
$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("D39").Value = "8.71"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +12.10%  "

$ws.Range("D41").Value = "408.64"
$ws.Range("E41").Value = "  +10.17%  "

$ws.Range("D42").Value = "0.0357"
$ws.Range("E42").Value = "  +3.95%  "

$ws.Range("D43").Value = "2.791.55"
$ws.Range("E43").Value = "  +4.69%  "

$ws.Range("D44").Value = "0.108"
$ws.Range("E44").Value = "  +0.76%  "

$ws.Range("E45").Value = "  +8.29%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").Value = "123.47"
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").Value = "34.54"
$ws.Range("E48").Value = "  +25.29%  "

$ws.Range("E49").Value = "  +4.40%  "

$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("D51").Value = "23.88"
$ws.Range("E51").Value = "  +3.42%  "

# Row 35/36: Filecoin and Mantle swap positions with updated values
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").Value = "  +6.73%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "6.00"
$ws.Range("E36").Value = "  +6.64%  "

